# Edit script: "Cambio en la función de reincoorporar empleado"
# Applies the textual changes described by the diff using Find/Replace.

$d = $word.ActiveDocument

function Replace-All($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2)
}

# 1. Header date
Replace-All "Oaxaca de Juárez, Oax. 30 de Octubre de 2025" "Oaxaca de Juárez, Oax. 3 de Noviembre de 2025"

# 2. Employee first name (drop "NOÉ") - 3 occurrences
Replace-All "ISIDRO NOÉ" "ISIDRO"

# 3. Age textbox number 20 -> 30
Replace-All "20" "30"

# 4. Department name - 2 occurrences
Replace-All "DEPARTAMENTO DE RECURSOS HUMANOS" "DEPARTAMENTO DE RECURSOS MATERIALES"

# 5. Date stamp
Replace-All "01 DE ENERO DE 2027" "05 DE NOVIEMBRE DE 2025"

# 6. ID / phone-like number - 3 occurrences
Replace-All "84849494949" "85493939399"

# 7. Municipality name -> Oaxaca (only the standalone occurrences, not the address line)
Replace-All "ZIMATLÁN DE ÁLVAREZ" "OAXACA"
